$d = $word.ActiveDocument

# Paragraph 1: headline block has two runs joined by a <w:br/> (date line + title line).
# Replace each run's text via Find/Replace so the line break between them is preserved.
$d.Content.Find.Execute('⚡️🚀המאמר היומי של מייק -29.11.24: ⚡️🚀', $true, $false, $false, $false, $false, $true, 1, $false, '⚡️🚀המאמר היומי של מייק -28.11.24: ⚡️🚀', 2) | Out-Null
$d.Content.Find.Execute('In-Context Learning with Long-Context Models: An In-Depth Exploration', $true, $false, $false, $false, $false, $true, 1, $false, 'Parameter-Efficient Fine-Tuning with Discrete Fourier Transform', 2) | Out-Null

# Paragraphs 2-9: swap the old review body for the new one (same paragraph count / positions).
$d.Content.Find.Execute('המאמר מציג מחקר אמפירי מקיף של למידה in-context או ICL עם מודלי שפה בעלי חלון הקשר ארוך. אזכיר שעם ICL המודל מקבל כמה דוגמאות המדגימות פעולות מסוימות ולאחר מכן המודל מתבקש לבצע פעולה זו על דוגמאות חדשות.', $true, $false, $false, $false, $false, $true, 1, $false, 'רקע: PeFT:', 2) | Out-Null
$d.Content.Find.Execute('ממצאים חדשים על התנהגות של ICL ל-LLMs בעלי חלון הקשר ארוך:', $true, $false, $false, $false, $false, $true, 1, $false, 'נתחיל את הסקירה ברענון קצרצר לגבי שיטות טיוב (fine-tuning) חסכוניות של מודלי שפה. PeFT הינה משפחה של שיטות המאפשרות טיוב של מודלים גדולים (בפרט מודל שפה) תוך שימוש במספר מצומצם של פרמטרים, מה שחוסך משמעותית במשאבי חישוב וזיכרון. ', 2) | Out-Null
$d.Content.Find.Execute('  1. שיפור ביצועים מתמשך: עלייה משמעותית בביצועים כאשר מעלים את מספר הדוגמאות בהדגמה מ-10 ל-1000 דוגמאות', $true, $false, $false, $false, $false, $true, 1, $false, 'רקע: LoRA:', 2) | Out-Null
$d.Content.Find.Execute('  2. רגישות פחותה לסדר: השפעת סדר הדוגמאות יורדת ב-50% ב-1000 דוגמאות לעומת 10(עבור סידור אקראי)', $true, $false, $false, $false, $false, $true, 1, $false, 'אחת השיטות הפופולריות ביותר ב-PeFT, הנקראת LoRA, מקפיאה את משקולות המודל ומאמנת מטריצות תוספת לכל שכבה של הטרנספורמטורים. כל מטריצת תוספת נלמדת הינה בעלת בדרגה נמוכה (low-rank), כך שניתן לייצגה על ידי מכפלה של שתי מטריצות קטנות (במימד האמצעי של המכפלה). ', 2) | Out-Null
$d.Content.Find.Execute('  3. ירידה ביתרון ה-RAG: היתרון של RAG פוחת משמעותית עם יותר דוגמאות', $true, $false, $false, $false, $false, $true, 1, $false, 'היתרון המרכזי של LoRA הוא שהיא מאפשרת להתאים מודלים גדולים למשימות ספציפיות תוך אימון של חלק קטן (נגיד 1% מכלל הפרמטרים שלו), מה שהופך אותה ליעילה במיוחד. שיטה זו הוכיחה את עצמה כאפקטיבית במיוחד בהתאמת מודלי שפה גדולים למשימות ספציפיות. בנוסף, LoRA מאפשרת החלפה מהירה בין גרסאות שונות של המודל המטויב, מכיוון שניתן לשמור את המטריצות הקטנות בנפרד מהמודל המקורי.', 2) | Out-Null
$d.Content.Find.Execute('  4. השפעת קיבוץ דוגמאות לפי קטגוריות: מיון דוגמאות לפי קטגוריות פוגע יותר בביצועים ככל שחלון ההקשר גדל', $true, $false, $false, $false, $false, $true, 1, $false, 'שיטה מוצעת:', 2) | Out-Null
$d.Content.Find.Execute('  5. יעילות אורכי attention קצרים: ניתן להשיג ביצועים דומים עם מנגנון attention קצר יחסית המשתרע ל-50-75 דוגמאות', $true, $false, $false, $false, $false, $true, 1, $false, 'הרעיון המרכזי הוא להסתכל על שינויי המשקולות של רשת הנוירונים כמו על תמונה או אות, ולייצג אותם בציר התדר במקום ערכים ישירים. כשאנחנו רוצים לטייב את המודל, במקום לשנות את כל המשקולות באופן ישיר (שדורש המון פרמטרים), אנחנו:', 2) | Out-Null
$d.Content.Find.Execute('  6. השוואה לטיוב (fine-tuning): למידת in-context לאורכי חלון הקשר ארוכים לרוב משתווה או עולה על טיוב עם מעט דוגמאות אולם הטיוב מנצח כאשר יש מספיק דוגמאות.', $true, $false, $false, $false, $false, $true, 1, $false, '1. מגדירים מראש כמה נקודות דגימה במרחב התדרים שבהן נרצה להתמקד. זה כמו לבחור אילו תדרים אנחנו רוצים לשמור בייצוג הדחוס שלנו. זה נעשה על ידי בחירת מטריצת תדרים קבועה (לא נלמדת) E בגודל 2xn המשמשת לבניית ייצוג של מטריצת תוספת. מטריצה זו היא קבועה לכל השכבות של הטרנספורמרים.', 2) | Out-Null

# New paragraphs describing steps 2-4, advantages, "it works because", and the summary,
# inserted right after the "step 1" paragraph and before the closing arXiv link paragraph.
$anchor = $d.Paragraphs(9).Range
$anchor.InsertParagraphAfter()
$d.Paragraphs(10).Range.Text = '2. לומדים וקטור c בגודל n (לכל שכבה) כאשר דרך שילובו עם E בונים את מטריצת התוספות בתחום התדר F (הסבר לאיך זה נבנה לא נראה ברור במאמר)'
$d.Paragraphs(10).Range.InsertParagraphAfter()
$d.Paragraphs(11).Range.Text = '3. מעבירים את F דרך Gaussian bandpass filter (כלומר דוגמים בעיקר תדרים נמוכים, הנמצאים קרוב למרכז המטריצה).'
$d.Paragraphs(11).Range.InsertParagraphAfter()
$d.Paragraphs(12).Range.Text = '4. מעבירים את מטריצת F לתחום הזמן (הרגיל) ומשתמשים בה בדיוק כמו ב-LoRA'
$d.Paragraphs(12).Range.InsertParagraphAfter()
$d.Paragraphs(13).Range.Text = 'יתרונות השיטה המוצעת:'
$d.Paragraphs(13).Range.InsertParagraphAfter()
$d.Paragraphs(14).Range.Text = 'היתרון הגדול הוא שתדרים הם דרך מאוד יעילה לייצג מידע (צריך 2n+ Ln משקלים כאשר L מספר השכבות במודל). בדיוק כמו שאפשר לדחוס תמונה או מוזיקה על ידי שמירת התדרים החשובים ביותר, כאן אנחנו יכולים לייצג שינויים מורכבים במשקולות באמצעות מספר קטן מאוד של תדרים.'
$d.Paragraphs(14).Range.InsertParagraphAfter()
$d.Paragraphs(15).Range.Text = 'זה עובד טוב(כנראה):'
$d.Paragraphs(15).Range.InsertParagraphAfter()
$d.Paragraphs(16).Range.Text = '- שינויים במשקולות נוטים להיות "חלקים" יחסית, כלומר יש בהם מבנה שאפשר לתפוס טוב עם תדרים'
$d.Paragraphs(16).Range.InsertParagraphAfter()
$d.Paragraphs(17).Range.Text = '- הבסיס המתמטי של פורייה הוא אורתוגונלי, מה שאומר שכל תדר מוסיף מידע ייחודי'
$d.Paragraphs(17).Range.InsertParagraphAfter()
$d.Paragraphs(18).Range.Text = '- אנחנו יכולים לבחור מראש כמה תדרים אנחנו רוצים לשמור, ובכך לשלוט ישירות בכמות הפרמטרים'
$d.Paragraphs(18).Range.InsertParagraphAfter()
$d.Paragraphs(19).Range.Text = 'סיכום:'
$d.Paragraphs(19).Range.InsertParagraphAfter()
$d.Paragraphs(20).Range.Text = 'בניגוד לשיטות אחרות שמנסות להקטין את כמות הפרמטרים על ידי הגבלת הדרגה של המטריצות (כמו LoRA), הגישה הזו מסתכלת על הבעיה מזווית שונה - דרך עדשת התדרים, ומצליחה להשיג דחיסה משמעותית יותר.'

# Final paragraph: update the arXiv link to the new paper.
$d.Content.Find.Execute('https://arxiv.org/abs/2405.00200', $true, $false, $false, $false, $false, $true, 1, $false, 'https://arxiv.org/abs/2405.03003', 2) | Out-Null

